$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Business demography" row (row 10) with the latest period data
$ws.Range("C10").Value = "Dec 2022 - Dec 2023 (18/11/24)"
$ws.Range("D10").Value = "Dec 2023 - Dec 2024 (Nov 25)"

# Move the active selection to D11 (matches saved cursor position in the source file)
$ws.Range("D11").Select()
